# Atualizacao de bases das ligas, do dia: 14-04-2024 as 15:12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 273 is brand new; copy formatting for the styled columns (A: bold/border,
# E: date format) from the row above before writing values into it.
$ws.Cells.Item(272, 1).Copy($ws.Cells.Item(273, 1))
$ws.Cells.Item(272, 5).Copy($ws.Cells.Item(273, 5))

    # Row 235
    $ws.Cells.Item(235, 2).Value = 6836277
    $ws.Cells.Item(235, 6).Value = "CFR Cluj"
    $ws.Cells.Item(235, 7).Value = "AFC Hermannstadt"
    $ws.Cells.Item(235, 11).Value = 1.7
    $ws.Cells.Item(235, 12).Value = 3.4
    $ws.Cells.Item(235, 13).Value = 5
    $ws.Cells.Item(235, 14).Value = 1.65
    $ws.Cells.Item(235, 15).Value = 3.5
    $ws.Cells.Item(235, 16).Value = 5.25
    $ws.Cells.Item(235, 17).Value = -0.75
    $ws.Cells.Item(235, 18).Value = 1.85
    $ws.Cells.Item(235, 19).Value = 2
    $ws.Cells.Item(235, 21).Value = 1.875
    $ws.Cells.Item(235, 22).Value = 1.975
    $ws.Cells.Item(235, 23).Value = 0.6499999999999999
    $ws.Cells.Item(235, 26).Value = 0.425
    $ws.Cells.Item(235, 27).Value = -0.5
    $ws.Cells.Item(235, 29).Value = 0.9750000000000001
    # Row 236
    $ws.Cells.Item(236, 2).Value = 6870268
    $ws.Cells.Item(236, 6).Value = "Petrolul Ploiesti"
    $ws.Cells.Item(236, 7).Value = "ACS Sepsi"
    $ws.Cells.Item(236, 8).Value = 1
    $ws.Cells.Item(236, 9).Value = 2
    $ws.Cells.Item(236, 10).Value = "A"
    $ws.Cells.Item(236, 11).Value = 2.8
    $ws.Cells.Item(236, 12).Value = 3
    $ws.Cells.Item(236, 13).Value = 2.55
    $ws.Cells.Item(236, 14).Value = 3
    $ws.Cells.Item(236, 15).Value = 3.2
    $ws.Cells.Item(236, 16).Value = 2.3
    $ws.Cells.Item(236, 18).Value = 1.85
    $ws.Cells.Item(236, 19).Value = 2
    $ws.Cells.Item(236, 21).Value = 1.875
    $ws.Cells.Item(236, 22).Value = 1.975
    $ws.Cells.Item(236, 24).Value = -1
    $ws.Cells.Item(236, 25).Value = 1.3
    $ws.Cells.Item(236, 26).Value = -1
    $ws.Cells.Item(236, 27).Value = 1
    $ws.Cells.Item(236, 28).Value = 0.875
    $ws.Cells.Item(236, 29).Value = -1
    # Row 237
    $ws.Cells.Item(237, 2).Value = 6865915
    $ws.Cells.Item(237, 6).Value = "FC Voluntari"
    $ws.Cells.Item(237, 7).Value = "Universitatea Cluj"
    $ws.Cells.Item(237, 8).Value = 0
    $ws.Cells.Item(237, 9).Value = 0
    $ws.Cells.Item(237, 10).Value = "D"
    $ws.Cells.Item(237, 11).Value = 3.5
    $ws.Cells.Item(237, 12).Value = 3.25
    $ws.Cells.Item(237, 13).Value = 2.05
    $ws.Cells.Item(237, 14).Value = 3.4
    $ws.Cells.Item(237, 15).Value = 3.1
    $ws.Cells.Item(237, 16).Value = 2.15
    $ws.Cells.Item(237, 18).Value = 1.975
    $ws.Cells.Item(237, 19).Value = 1.875
    $ws.Cells.Item(237, 21).Value = 2.05
    $ws.Cells.Item(237, 22).Value = 1.75
    $ws.Cells.Item(237, 24).Value = 2.1
    $ws.Cells.Item(237, 25).Value = -1
    $ws.Cells.Item(237, 26).Value = 0.4875
    $ws.Cells.Item(237, 27).Value = -0.5
    $ws.Cells.Item(237, 28).Value = -1
    $ws.Cells.Item(237, 29).Value = 0.75
    # Row 239
    $ws.Cells.Item(239, 2).Value = 6852370
    $ws.Cells.Item(239, 6).Value = "Dinamo Bucharest"
    $ws.Cells.Item(239, 7).Value = "ACS UTA Batrana Doamna"
    $ws.Cells.Item(239, 11).Value = 2.55
    $ws.Cells.Item(239, 12).Value = 2.875
    $ws.Cells.Item(239, 13).Value = 3
    $ws.Cells.Item(239, 14).Value = 2.375
    $ws.Cells.Item(239, 15).Value = 3
    $ws.Cells.Item(239, 16).Value = 3.1
    $ws.Cells.Item(239, 17).Value = -0.25
    $ws.Cells.Item(239, 18).Value = 2
    $ws.Cells.Item(239, 19).Value = 1.85
    $ws.Cells.Item(239, 21).Value = 1.975
    $ws.Cells.Item(239, 22).Value = 1.875
    $ws.Cells.Item(239, 23).Value = 1.375
    $ws.Cells.Item(239, 26).Value = 1
    $ws.Cells.Item(239, 27).Value = -1
    $ws.Cells.Item(239, 29).Value = 0.875
    # Row 267
    $ws.Cells.Item(267, 2).Value = 7951792
    $ws.Cells.Item(267, 5).Value = 45394.47916666666
    $ws.Cells.Item(267, 6).Value = "ACS UTA Batrana Doamna"
    $ws.Cells.Item(267, 7).Value = "AFC Hermannstadt"
    $ws.Cells.Item(267, 8).Value = 1
    $ws.Cells.Item(267, 9).Value = 3
    $ws.Cells.Item(267, 10).Value = "A"
    $ws.Cells.Item(267, 11).Value = 2.3
    $ws.Cells.Item(267, 13).Value = 3.1
    $ws.Cells.Item(267, 15).Value = 2.875
    $ws.Cells.Item(267, 16).Value = 2.7
    $ws.Cells.Item(267, 18).Value = 1.975
    $ws.Cells.Item(267, 19).Value = 1.875
    $ws.Cells.Item(267, 21).Value = 1.8
    $ws.Cells.Item(267, 22).Value = 2.05
    $ws.Cells.Item(267, 23).Value = -1
    $ws.Cells.Item(267, 24).Value = -1
    $ws.Cells.Item(267, 25).Value = 1.7
    $ws.Cells.Item(267, 26).Value = -1
    $ws.Cells.Item(267, 27).Value = 0.875
    $ws.Cells.Item(267, 28).Value = 0.8
    $ws.Cells.Item(267, 29).Value = -1
    # Row 268
    $ws.Cells.Item(268, 2).Value = 7951755
    $ws.Cells.Item(268, 5).Value = 45394.60416666666
    $ws.Cells.Item(268, 6).Value = "CS U Craiova"
    $ws.Cells.Item(268, 7).Value = "Farul Constanta"
    $ws.Cells.Item(268, 8).Value = 1
    $ws.Cells.Item(268, 9).Value = 2
    $ws.Cells.Item(268, 10).Value = "A"
    $ws.Cells.Item(268, 11).Value = 1.8
    $ws.Cells.Item(268, 12).Value = 3.4
    $ws.Cells.Item(268, 13).Value = 4.2
    $ws.Cells.Item(268, 14).Value = 1.8
    $ws.Cells.Item(268, 15).Value = 3.6
    $ws.Cells.Item(268, 16).Value = 4
    $ws.Cells.Item(268, 18).Value = 1.825
    $ws.Cells.Item(268, 19).Value = 2.025
    $ws.Cells.Item(268, 20).Value = 2.75
    $ws.Cells.Item(268, 21).Value = 1.95
    $ws.Cells.Item(268, 22).Value = 1.9
    $ws.Cells.Item(268, 23).Value = -1
    $ws.Cells.Item(268, 24).Value = -1
    $ws.Cells.Item(268, 25).Value = 3
    $ws.Cells.Item(268, 26).Value = -1
    $ws.Cells.Item(268, 27).Value = 1.025
    $ws.Cells.Item(268, 28).Value = 0.475
    $ws.Cells.Item(268, 29).Value = -0.5
    # Row 269
    $ws.Cells.Item(269, 2).Value = 7951791
    $ws.Cells.Item(269, 5).Value = 45395.64583333334
    $ws.Cells.Item(269, 6).Value = "Universitatea Cluj"
    $ws.Cells.Item(269, 7).Value = "Petrolul Ploiesti"
    $ws.Cells.Item(269, 11).Value = 1.909
    $ws.Cells.Item(269, 12).Value = 3.2
    $ws.Cells.Item(269, 13).Value = 4
    $ws.Cells.Item(269, 14).Value = 1.909
    $ws.Cells.Item(269, 15).Value = 3.2
    $ws.Cells.Item(269, 16).Value = 4.2
    $ws.Cells.Item(269, 17).Value = -0.5
    $ws.Cells.Item(269, 18).Value = 1.95
    $ws.Cells.Item(269, 19).Value = 1.9
    $ws.Cells.Item(269, 20).Value = 2
    # Row 270
    $ws.Cells.Item(270, 2).Value = 7951788
    $ws.Cells.Item(270, 5).Value = 45396.3125
    $ws.Cells.Item(270, 6).Value = "FC Botosani"
    $ws.Cells.Item(270, 7).Value = "FC U Craiova 1948"
    $ws.Cells.Item(270, 11).Value = 2.4
    $ws.Cells.Item(270, 12).Value = 3
    $ws.Cells.Item(270, 13).Value = 2.9
    $ws.Cells.Item(270, 14).Value = 2.375
    $ws.Cells.Item(270, 15).Value = 3
    $ws.Cells.Item(270, 16).Value = 3
    $ws.Cells.Item(270, 18).Value = 2.1
    $ws.Cells.Item(270, 19).Value = 1.775
    $ws.Cells.Item(270, 20).Value = 2.25
    $ws.Cells.Item(270, 21).Value = 1.875
    $ws.Cells.Item(270, 22).Value = 1.975
    # Row 271
    $ws.Cells.Item(271, 2).Value = 7951789
    $ws.Cells.Item(271, 5).Value = 45396.51041666666
    $ws.Cells.Item(271, 6).Value = "Dinamo Bucharest"
    $ws.Cells.Item(271, 7).Value = "CSM Politehnica Iasi"
    $ws.Cells.Item(271, 11).Value = 2.25
    $ws.Cells.Item(271, 13).Value = 3.1
    $ws.Cells.Item(271, 14).Value = 2.1
    $ws.Cells.Item(271, 15).Value = 3.2
    $ws.Cells.Item(271, 16).Value = 3.2
    $ws.Cells.Item(271, 17).Value = -0.25
    $ws.Cells.Item(271, 18).Value = 1.875
    $ws.Cells.Item(271, 19).Value = 1.975
    $ws.Cells.Item(271, 20).Value = 2.5
    $ws.Cells.Item(271, 21).Value = 2.15
    $ws.Cells.Item(271, 22).Value = 1.725
    # Row 272
    $ws.Cells.Item(272, 2).Value = 7951754
    $ws.Cells.Item(272, 5).Value = 45396.625
    $ws.Cells.Item(272, 6).Value = "CFR Cluj"
    $ws.Cells.Item(272, 7).Value = "FCSB"
    $ws.Cells.Item(272, 11).Value = 2.3
    $ws.Cells.Item(272, 12).Value = 3.1
    $ws.Cells.Item(272, 13).Value = 3
    $ws.Cells.Item(272, 14).Value = 2.5
    $ws.Cells.Item(272, 15).Value = 3.1
    $ws.Cells.Item(272, 16).Value = 2.8
    $ws.Cells.Item(272, 17).Value = 0
    $ws.Cells.Item(272, 18).Value = 1.775
    $ws.Cells.Item(272, 19).Value = 2.1
    $ws.Cells.Item(272, 20).Value = 2.25
    $ws.Cells.Item(272, 21).Value = 1.85
    $ws.Cells.Item(272, 22).Value = 2
    # Row 273
    $ws.Cells.Item(273, 1).Value = 271
    $ws.Cells.Item(273, 2).Value = 7951756
    $ws.Cells.Item(273, 3).Value = "Romania Liga I"
    $ws.Cells.Item(273, 4).Value = "Romania Liga I"
    $ws.Cells.Item(273, 5).Value = 45397.60416666666
    $ws.Cells.Item(273, 6).Value = "Rapid Bucuresti"
    $ws.Cells.Item(273, 7).Value = "ACS Sepsi"
    $ws.Cells.Item(273, 11).Value = 1.95
    $ws.Cells.Item(273, 12).Value = 3.25
    $ws.Cells.Item(273, 13).Value = 3.75
    $ws.Cells.Item(273, 14).Value = 1.8
    $ws.Cells.Item(273, 15).Value = 3.3
    $ws.Cells.Item(273, 16).Value = 4.2
    $ws.Cells.Item(273, 17).Value = -0.5
    $ws.Cells.Item(273, 18).Value = 1.85
    $ws.Cells.Item(273, 19).Value = 2
    $ws.Cells.Item(273, 20).Value = 2.5
    $ws.Cells.Item(273, 21).Value = 1.95
    $ws.Cells.Item(273, 22).Value = 1.9
    $ws.Cells.Item(273, 23).Value = 0
    $ws.Cells.Item(273, 24).Value = 0
    $ws.Cells.Item(273, 25).Value = 0
    $ws.Cells.Item(273, 26).Value = 0
    $ws.Cells.Item(273, 27).Value = 0
